$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B (descriptions -> short names) to reflect the updated spec.
# Column A (the punctuator symbols) stays the same.
$ws.Range("B1").Value = "Name"
$ws.Range("B2").Value = "Semicolon"
$ws.Range("B3").Value = "Comma"
$ws.Range("B4").Value = "Dot"
$ws.Range("B5").Value = "Colon"
$ws.Range("B6").Value = "Left curly bracket"
$ws.Range("B7").Value = "Right curly bracket"
$ws.Range("B8").Value = "Left parenthesis"
$ws.Range("B9").Value = "Right parenthesis"
$ws.Range("B10").Value = "Left square bracket"
$ws.Range("B11").Value = "Right square bracket"
$ws.Range("B12").Value = "String delimiter"

# Narrow column B to fit the new, shorter text.
$ws.Columns.Item(2).ColumnWidth = 18

# Move the active selection to reflect where the editor ended up.
$ws.Range("B12").Select()
